# Insert a new "Industry" column right before the existing "Mutual Fund"
# column (column C), shifting Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/
# MoM/QoQ one column to the right (C->D, D->E, E->F, F->G, G->H, H->I, I->J).
# This reproduces the motilal_portfolio_change_engine output which now
# includes an Industry classification for each holding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at C; everything from the old column C onward shifts
# right by one, which also pushes the sheet's dimension from I27 to J27.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Industry"

# Industry values for each holding (rows 2-27), in row order.
$industries = @(
  "Auto Components",
  "Construction",
  "Power",
  "Industrial Products",
  "Construction",
  "Insurance",
  "Construction",
  "Metals & Minerals Trading",
  "Transport Infrastructure",
  "Realty",
  "Construction",
  "Power",
  "Cement & Cement Products",
  "Chemicals & Petrochemicals",
  "Electrical Equipment",
  "Industrial Products",
  "Minerals & Mining",
  "Construction",
  "Power",
  "Power",
  "Construction",
  "Power",
  "Electrical Equipment",
  "IT - Software",
  "IT - Software",
  "Petroleum Products"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $industries[$i]
}
